# Adds a "Savings" column (J) and three new expenditure rows, moving the
# "Total" row down from row 9 to row 12. Mirrors the commit
# "Added a savings column in the output".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Turn the old "Total" row (row 9) into a normal data row.
# ------------------------------------------------------------------
$ws.Range("A9").Value = "31-8-2022"
$ws.Range("B9").Value = "Food"
$ws.Range("C9").Value = 55
$ws.Range("D9").Value = "Gpay"
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = ""

# ------------------------------------------------------------------
# 2. Insert two more expenditure rows (10 and 11), same style pattern
#    as the rest of the data rows (A:D style 2, E:G style 3, H:I style 4).
# ------------------------------------------------------------------
$ws.Range("A10").Value = "31-8-2022"
$ws.Range("B10").Value = "Cab"
$ws.Range("C10").Value = 444
$ws.Range("D10").Value = "Gpay"

$ws.Range("A11").Value = "31-8-2022"
$ws.Range("B11").Value = "Cab"
$ws.Range("C11").Value = 444
$ws.Range("D11").Value = "Gpay"

# ------------------------------------------------------------------
# 3. Re-create the "Total" row at row 12, with the SUM ranges widened
#    to cover the new rows (3:11 instead of 3:8).
# ------------------------------------------------------------------
$ws.Range("A12").Value = "Total"
$ws.Range("B12").Value = ""
$ws.Range("C12").Formula = "=SUM(C3:C11)"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Formula = "=SUM(F3:F11)"
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = ""
$ws.Range("I12").Formula = "=SUM(I3:I11)"

# Copy the cell styles (fill/alignment) down column by column so rows
# 9-12 match the same look as the rest of the table.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A9:A12").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B9:B12").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C9:C12").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D9:D12").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E9:E12").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F9:F12").PasteSpecial(-4122) | Out-Null
$ws.Range("G3").Copy() | Out-Null
$ws.Range("G9:G12").PasteSpecial(-4122) | Out-Null
$ws.Range("H3").Copy() | Out-Null
$ws.Range("H9:H12").PasteSpecial(-4122) | Out-Null
$ws.Range("I3").Copy() | Out-Null
$ws.Range("I9:I12").PasteSpecial(-4122) | Out-Null

# Restore the values/formulas that the paste-special above would have
# clobbered back to "Total"-row shape (paste only copied formatting, so
# values are untouched, but make sure formulas are still correct).
$ws.Range("C12").Formula = "=SUM(C3:C11)"
$ws.Range("F12").Formula = "=SUM(F3:F11)"
$ws.Range("I12").Formula = "=SUM(I3:I11)"

# ------------------------------------------------------------------
# 4. New "Savings" column J: centered, no fill for the header/data
#    rows (1-8), centered + yellow fill for the total-block rows (9-12).
# ------------------------------------------------------------------
$yellow = 6740479  # RGB(255, 217, 102) -> FFD966

1..8 | ForEach-Object {
    $cell = $ws.Range("J$_")
    $cell.Value = ""
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

9..11 | ForEach-Object {
    $cell = $ws.Range("J$_")
    $cell.Value = ""
    $cell.Interior.Color = $yellow
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

$ws.Range("J12").Formula = "=I12-F12-C12"
$ws.Range("J12").Interior.Color = $yellow
$ws.Range("J12").HorizontalAlignment = -4108
$ws.Range("J12").VerticalAlignment = -4108

"done"
